$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column I (9) width: 14 -> 10 ---
# Raw OOXML width = ColumnWidth + 0.8333333 (5/6), so to land exactly on 10 we
# need to request 9.1666666...
$ws.Columns.Item(9).ColumnWidth = 9.166666666666666

# --- G2: reorder recorder email list ---
$ws.Range("G2").Value = "Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg"

# --- Class Statistics block ---
$ws.Range("L6").Value = 4   # Recorded Sessions: 3 -> 4
$ws.Range("L7").Value = 0   # Missing Sessions: 1 -> 0

# --- G9: reorder recorder email list ---
$ws.Range("G9").Value = "Safa.hany@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"

# Percent-looking values must stay literal text (not get auto-converted to a
# numeric percentage, which would rewrite the cell's number format/style).
# Trick: compute the literal text via a string-literal formula in a scratch
# cell, then Copy/PasteSpecial *values* into the target - this carries over
# the text as a shared string rather than re-parsing "13.8%" as 0.138.
$scratch = $ws.Range("ZZ1")

$scratch.Formula = '="13.8%"'
$scratch.Copy()
$ws.Range("L9").PasteSpecial(-4163)   # xlPasteValues

$scratch.Formula = '="23.8%"'
$scratch.Copy()
$ws.Range("L10").PasteSpecial(-4163)

$scratch.Formula = '="13.8%"'
$scratch.Copy()
$ws.Range("R15").PasteSpecial(-4163)

$scratch.Formula = '="23.8%"'
$scratch.Copy()
$ws.Range("S15").PasteSpecial(-4163)

$scratch.Value = ""

# --- Group Statistics row 15 ---
$ws.Range("O15").Value = 4   # Recorded: 3 -> 4
$ws.Range("P15").Value = 0   # Missing: 1 -> 0

# --- Row 26 (PHARMACOLOGY session 1) flips from Not Recorded -> Recorded ---
# Re-use the same look (green "Recorded" fill/font) as the other recorded
# rows (e.g. row 2) by copying their format only.
$src = $ws.Range("A2:I2")
$dst = $ws.Range("A26:I26")
$src.Copy()
$dst.PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("G26").Value = "nancy.abdelshafy@med.asu.edu.eg"
$ws.Range("H26").Value = "102/251"
$ws.Range("I26").Value = "Recorded"

$excel.CutCopyMode = $false
